$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data" (sheet1): add a "name" column and two data rows, making
# sure the date-looking strings stay TEXT instead of being auto-parsed
# into date serials (this is the "date format corruption" bug fix).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cell D1 = "name", matching the bold/bordered header style.
$ws1.Range("D1").Value = "name"
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# Row 2
$ws1.Range("A2").Value = "80b1787c-dc54-40f8-9be1-69c6deba7659"
$ws1.Range("A1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Range("B2").Value = "'12/10/2022"
$ws1.Range("B2").ClearFormats()

$ws1.Range("C2").Value = " "

$ws1.Range("D2").Value = "add all housing locs"

# Row 3
$ws1.Range("A3").Value = "1ee21dad-c8ff-49da-b74d-3af27f532e6c"
$ws1.Range("A1").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

$ws1.Range("B3").Value = "'12/10/2022"
$ws1.Range("B3").ClearFormats()

$ws1.Range("C3").Value = " "

$ws1.Range("D3").Value = "add jobs to objects as experiences"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Sheet "headers" (sheet2): convert the text "False" placeholders into
# real booleans, and add a "name" field description row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D2").Value = $False
$ws2.Range("E2").Value = $False
$ws2.Range("F2").Value = $False

$ws2.Range("D3").Value = $False
$ws2.Range("E3").Value = $False
$ws2.Range("F3").Value = $False

$ws2.Range("D4").Value = $False
$ws2.Range("E4").Value = $False
$ws2.Range("F4").Value = $False

# Row 5 - the new "name" field
$ws2.Range("A5").Value = "name"
$ws2.Range("A1").Copy()
$ws2.Range("A5").PasteSpecial(-4122)

$ws2.Range("B5").Value = "str"
$ws2.Range("C5").Value = "name"
$ws2.Range("D5").Value = $True
$ws2.Range("E5").Value = $True
$ws2.Range("F5").Value = $True

$excel.CutCopyMode = 0
